$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 2 (pushes existing rows 2-22 down to 3-23).
# Only the used columns (A:F) are shifted to avoid ballooning the sheet's
# used range.
$ws.Range("A2:F2").Insert(-4121)  # xlShiftDown

# Copy formatting (styles) from the row below (now row 3, the former row 2)
# onto the freshly inserted row 2 so borders/alignment/number-format match.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 2 with the latest Nalco circular entry.
$ws.Cells.Item(2, 1).Value = 22
$ws.Cells.Item(2, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(2, 3).Value = "IE07"
$ws.Cells.Item(2, 4).Value = 296.05
$ws.Cells.Item(2, 5).Value = "05-12-2025"
$ws.Cells.Item(2, 6).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-05-12-2025.pdf"

# Add the hyperlink for the new row's Circular Link cell.
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-05-12-2025.pdf")
